$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2406.125
$ws.Range("I43").Value = 1999.6666
$ws.Range("K43").Value = 1999.6666
$ws.Range("M43").Value = -1930.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 861.94446
$ws.Range("J112").Value = 836.94116
$ws.Range("L112").Value = 2510.82348
$ws.Range("N112").Value = -4726.82348

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 97242.57000000001
$ws.Range("J116").Value = 29749.5
$ws.Range("L116").Value = 29749.5
$ws.Range("N116").Value = -36633.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 81748.5
$ws.Range("J133").Value = 81748.5
$ws.Range("L133").Value = 81748.5
$ws.Range("N133").Value = -91868.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 55557972
$ws.Range("I137").Value = 111112380
$ws.Range("J137").Value = 3558.889
$ws.Range("K137").Value = 333337140
$ws.Range("L137").Value = 10676.667
$ws.Range("M137").Value = -333334590
$ws.Range("N137").Value = -15776.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2727.8235
$ws.Range("I138").Value = 1232.0555
$ws.Range("K138").Value = 3696.1665
$ws.Range("M138").Value = 1443.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1105.2307
$ws.Range("I32").Value = 966.25
$ws.Range("K32").Value = 966.25
$ws.Range("M32").Value = -679.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2260.4375
$ws.Range("I33").Value = 2011.9286
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 2011.9286
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -1682.9286
$ws.Range("N33").Value = -4658

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1095
$ws.Range("I61").Value = 884.6909000000001
$ws.Range("K61").Value = 884.6909000000001
$ws.Range("M61").Value = -672.6909000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1622.1
$ws.Range("I74").Value = 1185.1177
$ws.Range("K74").Value = 1185.1177
$ws.Range("M74").Value = -311.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1622.1
$ws.Range("I77").Value = 1185.1177
$ws.Range("K77").Value = 5925.5885
$ws.Range("M77").Value = -1557.5885

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1095
$ws.Range("I136").Value = 884.6909000000001
$ws.Range("K136").Value = 2654.0727
$ws.Range("M136").Value = -104.0727000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 91332.664
$ws.Range("J139").Value = 91332.664
$ws.Range("L139").Value = 91332.664
$ws.Range("N139").Value = -101612.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2034.4736
$ws.Range("I20").Value = 1995.8182
$ws.Range("J20").Value = 2087.625
$ws.Range("K20").Value = 1995.8182
$ws.Range("L20").Value = 2087.625
$ws.Range("M20").Value = -1748.8182
$ws.Range("N20").Value = -2581.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 38035.125
$ws.Range("I38").Value = 38468.715
$ws.Range("J38").Value = 35000
$ws.Range("K38").Value = 38468.715
$ws.Range("L38").Value = 35000
$ws.Range("M38").Value = -38052.715
$ws.Range("N38").Value = -35832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1821.7627
$ws.Range("I134").Value = 1149.7021
$ws.Range("K134").Value = 3449.1063
$ws.Range("M134").Value = -914.1062999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 833.3333
$ws.Range("I39").Value = 833.3333
$ws.Range("K39").Value = 833.3333
$ws.Range("M39").Value = -442.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 833.3333
$ws.Range("I49").Value = 833.3333
$ws.Range("K49").Value = 833.3333
$ws.Range("M49").Value = -651.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 80001070
$ws.Range("I132").Value = 86957420
$ws.Range("K132").Value = 260872260
$ws.Range("M132").Value = -260869730

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 396.85715
$ws.Range("I14").Value = 396.85715
$ws.Range("K14").Value = 1190.57145
$ws.Range("M14").Value = -1017.57145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1018.6667
$ws.Range("I34").Value = 202.18182
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 606.5454599999999
$ws.Range("L34").Value = 30000
$ws.Range("M34").Value = -522.5454599999999
$ws.Range("N34").Value = -30168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 97796
$ws.Range("J37").Value = 97796
$ws.Range("L37").Value = 293388
$ws.Range("N37").Value = -293612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2566.2
$ws.Range("I76").Value = 1209
$ws.Range("K76").Value = 3627
$ws.Range("M76").Value = -3244

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 2566.2
$ws.Range("I79").Value = 1209
$ws.Range("K79").Value = 3627
$ws.Range("M79").Value = -2301

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 18747.25
$ws.Range("I87").Value = 9971.333000000001
$ws.Range("J87").Value = 24012.8
$ws.Range("K87").Value = 29913.999
$ws.Range("L87").Value = 72038.39999999999
$ws.Range("M87").Value = -28665.999
$ws.Range("N87").Value = -74534.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 18747.25
$ws.Range("I90").Value = 9971.333000000001
$ws.Range("J90").Value = 24012.8
$ws.Range("K90").Value = 89741.997
$ws.Range("L90").Value = 216115.2
$ws.Range("M90").Value = -83501.997
$ws.Range("N90").Value = -228595.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 22500
$ws.Range("J54").Value = 22500
$ws.Range("L54").Value = 22500
$ws.Range("N54").Value = -23280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 77952
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 77952
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6511.2354
$ws.Range("J70").Value = 6926.3335
$ws.Range("L70").Value = 6926.3335
$ws.Range("N70").Value = -7466.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6511.2354
$ws.Range("J73").Value = 6926.3335
$ws.Range("L73").Value = 6926.3335
$ws.Range("N73").Value = -8798.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 25659296
$ws.Range("I132").Value = 37051040
$ws.Range("J132").Value = 27874.416
$ws.Range("K132").Value = 111153120
$ws.Range("L132").Value = 83623.24800000001
$ws.Range("M132").Value = -111150590
$ws.Range("N132").Value = -88683.24800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3674.1333
$ws.Range("I68").Value = 2610.1
$ws.Range("K68").Value = 2610.1
$ws.Range("M68").Value = -1861.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3674.1333
$ws.Range("I71").Value = 2610.1
$ws.Range("K71").Value = 13050.5
$ws.Range("M71").Value = -9306.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6374.25
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3713
$ws.Range("I132").Value = 3443.2222
$ws.Range("K132").Value = 10329.6666
$ws.Range("M132").Value = -7799.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5571.143
$ws.Range("I136").Value = 3999.6667
$ws.Range("K136").Value = 11999.0001
$ws.Range("M136").Value = -9449.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 40000
$ws.Range("I56").Value = 40000
$ws.Range("K56").Value = 40000
$ws.Range("M56").Value = -39286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1639
$ws.Range("I113").Value = 500.7647
$ws.Range("K113").Value = 1502.2941
$ws.Range("M113").Value = 667.7058999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -5800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4696926
$ws.Range("I136").Value = 4903834.5
$ws.Range("K136").Value = 14711503.5
$ws.Range("M136").Value = -14708953.5
